# Apply "more work towards final product" edit to 17_10_stimuli sheet.
#
# Adds a "carrier" value (column D) to the practice-pair rows (2-5) and
# fills in D for each generic-pair row's practice carrier; also tags the
# unique-video / unique-audio pair_kind (column J) for pairs E-H (rows 6-9)
# and adds four new rows (14-21) describing the unique_video / unique_audio
# carrier combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows: carrier column (D) gets the carrier word used in that
# practice trial.
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic pair rows E-H (6-9): tag pair_kind (J) as unique_video / unique_audio.
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# New rows 14-21: kind (C) + carrier (D) for the unique_video / unique_audio
# stimuli numbers 9-16.
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
